$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet: By_Odds_Bin
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("By_Odds_Bin")

# Row 2: (0, 5]
$ws1.Range("B2").Value = 3
$ws1.Range("C2").Value = -3
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = -3
$ws1.Range("F2").Value = 0

# Row 3: (5, 10]
$ws1.Range("B3").Value = 100
$ws1.Range("C3").Value = 56
$ws1.Range("D3").Value = 136
$ws1.Range("E3").Value = -80
$ws1.Range("F3").Value = 20

# Row 4: (10, 15]
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 0
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = ""

# -----------------------------------------------------------------
# Sheet: By_Field_Size
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("By_Field_Size")

# Row 2: 1-4
$ws2.Range("B2").Value = 6
$ws2.Range("C2").Value = 0.5
$ws2.Range("D2").Value = 5.5
$ws2.Range("E2").Value = -5
$ws2.Range("F2").Value = 16.7

# Row 3: 5
$ws2.Range("B3").Value = 19
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 17
$ws2.Range("E3").Value = -17
$ws2.Range("F3").Value = 10.5

# Row 4: 6
$ws2.Range("B4").Value = 36
$ws2.Range("C4").Value = 15
$ws2.Range("D4").Value = 44
$ws2.Range("E4").Value = -29
$ws2.Range("F4").Value = 19.4

# Row 5: 7
$ws2.Range("B5").Value = 42
$ws2.Range("C5").Value = 37.5
$ws2.Range("D5").Value = 69.5
$ws2.Range("E5").Value = -32
$ws2.Range("F5").Value = 23.8

# -----------------------------------------------------------------
# Sheet: By_Track
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("By_Track")

# Row 2: HAMILTON -> NEWMARKET
$ws3.Range("A2").Value = "NEWMARKET"
$ws3.Range("B2").Value = 103
$ws3.Range("C2").Value = 53
$ws3.Range("D2").Value = 136
$ws3.Range("E2").Value = -83
$ws3.Range("F2").Value = 19.4
